$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.17212692718101152
$ws.Range("A2").Value = -0.041489237629498632
$ws.Range("A3").Value = -0.0039999999770863326
$ws.Range("A4").Value = -0.0079999999575353087
$ws.Range("A5").Value = 0.043977678806101039
$ws.Range("A6").Value = -0.0019999999768156584
$ws.Range("A7").Value = -0.0099999999403852513
$ws.Range("A8").Value = -0.009999999938213211
$ws.Range("A9").Value = -0.0019999999719910733
$ws.Range("A10").Value = -0.0019999999699802373
$ws.Range("A11").Value = -0.0029999999651391107
$ws.Range("A12").Value = -0.0034999999622691291
$ws.Range("A13").Value = -0.0034999999601970089
$ws.Range("A14").Value = -0.0079999999394422261
$ws.Range("A15").Value = -0.00099999997067179436
$ws.Range("A16").Value = -0.0019999999658040224
$ws.Range("A17").Value = -0.0019999999654158884
$ws.Range("A18").Value = -0.003999999956245226
$ws.Range("A19").Value = -0.0039999999812523335
$ws.Range("A20").Value = -0.0039999999799213981
$ws.Range("A21").Value = -0.0039999999796958008
$ws.Range("A22").Value = -0.0039999999795261587
$ws.Range("A23").Value = -0.0049999999710710341
$ws.Range("A24").Value = -0.019999999900790932
$ws.Range("A25").Value = -0.019999999899435572
$ws.Range("A26").Value = 0.02693506368354015
$ws.Range("A27").Value = -0.0024999999724948374
$ws.Range("A28").Value = -0.0019999999689055414
$ws.Range("A29").Value = -0.0069999999422130088
$ws.Range("A30").Value = -0.059999999701431328
$ws.Range("A31").Value = 0.0522356490277609
$ws.Range("A32").Value = -0.0099999999253181926
$ws.Range("A33").Value = -0.039507653281411592
